$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.079.72'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.87'
$ws.Range("E3").Value = '  +2.07%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.73'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("E6").Value = '  +2.75%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.99'
$ws.Range("E8").Value = '  +6.86%  '
$ws.Range("E9").Value = '  +0.72%  '
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0983'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.115.43'
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("E13").Value = '  +3.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.854.52'
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.676'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.077.14'
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.02'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.00'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.16'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("E22").Value = '  +2.00%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +3.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.24'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.95'
$ws.Range("E26").Value = '  +3.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.56'
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("E28").Value = '  +3.81%  '
$ws.Range("E29").Value = '  +11.64%  '
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.64'
$ws.Range("E34").Value = '  +23.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.97'
$ws.Range("E35").Value = '  +10.49%  '
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("E38").Value = '  +11.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '90.62'
$ws.Range("E40").Value = '  +4.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.347.11'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.31'
$ws.Range("E43").Value = '  +4.44%  '
$ws.Range("E44").Value = '  +3.45%  '
$ws.Range("E45").Value = '  -3.68%  '
$ws.Range("B46").Value = 'Gas'
$ws.Range("C46").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.08'
$ws.Range("E46").Value = '  +75.10%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0532'
$ws.Range("E47").Value = '  +3.71%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.36'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.028.16'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("E50").Value = '  +16.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0674'
$ws.Range("E51").Value = '  +0.85%  '
